$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 99, shifting existing rows 99:144 down to 100:145.
$ws.Rows.Item(99).Insert()

# Populate the newly inserted row 99 with the new data record.
$ws.Range("A99").Value = 3
$ws.Range("B99").Value = "Femacal de La Calera"
$ws.Range("C99").Value = "Coquimbo"
$ws.Range("D99").Value = 44609
$ws.Range("E99").Value = 5
$ws.Range("F99").Value = 100112030
$ws.Range("G99").Value = "Poroto granado"
$ws.Range("H99").Value = "Sin especificar"
$ws.Range("I99").Value = "Primera"
$ws.Range("J99").Value = 70
$ws.Range("K99").Value = 25000
$ws.Range("L99").Value = 26000
$ws.Range("M99").Value = 25500
$ws.Range("N99").Value = "$/malla 25 kilos"
$ws.Range("O99").Value = "Provincia de Quillota"
$ws.Range("P99").Value = 1020
$ws.Range("Q99").Value = 25
$ws.Range("R99").Value = "Hortaliza"
